$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-17 (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then the 16 numeric NATMI metric columns E..T).
# A new "Resolving-Mac" sending cluster was added (previously it only existed as a
# target cluster), expanding the grid from 3x4=12 data rows to 4x4=16 data rows,
# and all TPM-derived metrics were recalculated with the new data.
$data = @(
    @("ECs", "Inhbb", "Acvr2b", "ECs", 3, 1, 1.165518666666667, 3.496556, 0.1888647065994748, 0.1888647065994748, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 2.000795777764, 18.007161999876, 0.07254078144577611, 0.07254078144577611),
    @("ECs", "Inhbb", "Acvr2b", "FAPs", 3, 1, 1.165518666666667, 3.496556, 0.1888647065994748, 0.1888647065994748, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.4099719680124445, 3.689747712112, 0.01486392927304215, 0.01486392927304215),
    @("ECs", "Inhbb", "Acvr2b", "MuSCs", 3, 1, 1.165518666666667, 3.496556, 0.1888647065994748, 0.1888647065994748, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 1.955241869183555, 17.597176822652, 0.07088918053624806, 0.07088918053624806),
    @("ECs", "Inhbb", "Acvr2b", "Resolving-Mac", 3, 1, 1.165518666666667, 3.496556, 0.1888647065994748, 0.1888647065994748, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.8431940908937777, 7.588746818043999, 0.03057081534440849, 0.03057081534440849),
    @("FAPs", "Inhbb", "Acvr2b", "ECs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.06060032033836409, 0.06060032033836409, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 0.641987946013, 5.777891514117001, 0.02327589242246235, 0.02327589242246235),
    @("FAPs", "Inhbb", "Acvr2b", "FAPs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.06060032033836409, 0.06060032033836409, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.1315461900671111, 1.183915710604, 0.004769334046849632, 0.004769334046849631),
    @("FAPs", "Inhbb", "Acvr2b", "MuSCs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.06060032033836409, 0.06060032033836409, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 0.627371231739889, 5.646341085659, 0.02274594934315114, 0.02274594934315114),
    @("FAPs", "Inhbb", "Acvr2b", "Resolving-Mac", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.06060032033836409, 0.06060032033836409, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.2705525713914445, 2.434973142523, 0.009809144525900969, 0.009809144525900967),
    @("MuSCs", "Inhbb", "Acvr2b", "ECs", 3, 1, 4.624333, 13.872999, 0.7493430352008683, 0.7493430352008682, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 7.938393614781, 71.445542533029, 0.2878141200817235, 0.2878141200817234),
    @("MuSCs", "Inhbb", "Acvr2b", "FAPs", 3, 1, 4.624333, 13.872999, 0.7493430352008683, 0.7493430352008682, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 1.626612215638667, 14.639509940748, 0.05897439535960084, 0.05897439535960082),
    @("MuSCs", "Inhbb", "Acvr2b", "MuSCs", 3, 1, 4.624333, 13.872999, 0.7493430352008683, 0.7493430352008682, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 7.757653100920333, 69.818877908283, 0.2812611983592395, 0.2812611983592394),
    @("MuSCs", "Inhbb", "Acvr2b", "Resolving-Mac", 3, 1, 4.624333, 13.872999, 0.7493430352008683, 0.7493430352008682, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 3.345472167405667, 30.109249506651, 0.1212933214003047, 0.1212933214003046),
    @("Resolving-Mac", "Inhbb", "Acvr2b", "ECs", 1, 0.3333333333333333, 0.007355666666666667, 0.022067, 0.001191937861292829, 0.001191937861292829, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 0.012627156673, 0.113644410057, 0.0004578097488396986, 0.0004578097488396986),
    @("Resolving-Mac", "Inhbb", "Acvr2b", "FAPs", 1, 0.3333333333333333, 0.007355666666666667, 0.022067, 0.001191937861292829, 0.001191937861292829, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.002587360653777778, 0.023286245884, 0.00009380725698894028, 0.00009380725698894026),
    @("Resolving-Mac", "Inhbb", "Acvr2b", "MuSCs", 1, 0.3333333333333333, 0.007355666666666667, 0.022067, 0.001191937861292829, 0.001191937861292829, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 0.01233966289322222, 0.111056966039, 0.0004473863844575594, 0.0004473863844575593),
    @("Resolving-Mac", "Inhbb", "Acvr2b", "Resolving-Mac", 1, 0.3333333333333333, 0.007355666666666667, 0.022067, 0.001191937861292829, 0.001191937861292829, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.005321454598111111, 0.047893091383, 0.0001929344710066311, 0.0001929344710066311)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
